$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.996.32'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.690.81'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.10'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.00'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000207'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +8.87%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.07'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.171.95'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.799.43'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.694.68'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.68%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.87'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.72'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '357.02'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000113'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +16.68%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.67%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.65'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.26'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '533.47'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.63'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.47'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.69'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.41'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.78'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.35'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0634'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.64'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.75'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0996'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.40%  '
